# Insert a new data row at row 304 (pushing the existing rows 304-365 down
# to 305-366) and populate it with the new "Arveja Verde" observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(304).Insert()

$ws.Cells.Item(304, 1).Value  = 6
$ws.Cells.Item(304, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(304, 3).Value  = "Metropolitana"
$ws.Cells.Item(304, 4).Value  = 45244
$ws.Cells.Item(304, 5).Value  = 13
$ws.Cells.Item(304, 6).Value  = 100112022
$ws.Cells.Item(304, 7).Value  = "Arveja Verde"
$ws.Cells.Item(304, 8).Value  = "Sin especificar"
$ws.Cells.Item(304, 9).Value  = "Primera"
$ws.Cells.Item(304, 10).Value = 400
$ws.Cells.Item(304, 11).Value = 18000
$ws.Cells.Item(304, 12).Value = 20000
$ws.Cells.Item(304, 13).Value = 19000
$ws.Cells.Item(304, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(304, 15).Value = "Región del Maule"
$ws.Cells.Item(304, 16).Value = 760
$ws.Cells.Item(304, 17).Value = 25
$ws.Cells.Item(304, 18).Value = "Hortaliza"
